$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2-484).
# The commit updates that "last changed" date from 2023-09-21 (45190)
# to 2023-09-23 (45192) for every row, leaving all other cells/styles intact.
$ws.Range("C2:C484").Value = 45192
